$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 497.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 497.5
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 497.5
$ws.Range("N32").Value = -1149.5
$ws.Range("M32").ClearContents()
# Row 138
$ws.Range("H138").Value = 3151.718
$ws.Range("I138").Value = 1516.68
$ws.Range("J138").Value = 6071.4287
$ws.Range("K138").Value = 4550.04
$ws.Range("L138").Value = 18214.2861
$ws.Range("M138").Value = 589.96
$ws.Range("N138").Value = -28494.2861
# Row 141
$ws.Range("H141").Value = 1223031.1
$ws.Range("I141").Value = 2012.7222
$ws.Range("J141").Value = 8549141
$ws.Range("K141").Value = 6038.1666
$ws.Range("L141").Value = 25647423
$ws.Range("M141").Value = -858.1665999999996
$ws.Range("N141").Value = -25657783

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5075.566
$ws.Range("I32").Value = 3771.848
$ws.Range("J32").Value = 13642.857
$ws.Range("K32").Value = 3771.848
$ws.Range("L32").Value = 13642.857
$ws.Range("M32").Value = -3484.848
$ws.Range("N32").Value = -14216.857
# Row 132
$ws.Range("H132").Value = 28575842
$ws.Range("I132").Value = 43482370
$ws.Range("K132").Value = 130447110
$ws.Range("M132").Value = -130444580
# Row 134
$ws.Range("H134").Value = 26952.666
$ws.Range("J134").Value = 26952.666
$ws.Range("L134").Value = 26952.666
$ws.Range("N134").Value = -37092.666
# Row 135
$ws.Range("H135").Value = 32485.715
$ws.Range("J135").Value = 32485.715
$ws.Range("L135").Value = 32485.715
$ws.Range("N135").Value = -42625.715
# Row 139
$ws.Range("H139").Value = 29750
$ws.Range("J139").Value = 29750
$ws.Range("L139").Value = 29750
$ws.Range("N139").Value = -40030
# Row 141
$ws.Range("H141").Value = 68791.125
$ws.Range("J141").Value = 68791.125
$ws.Range("L141").Value = 68791.125
$ws.Range("N141").Value = -79151.125

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 19233526
$ws.Range("I58").Value = 1398.2142
$ws.Range("J58").Value = 41671010
$ws.Range("K58").Value = 1398.2142
$ws.Range("L58").Value = 41671010
$ws.Range("M58").Value = -1195.2142
$ws.Range("N58").Value = -41671416
# Row 99
$ws.Range("H99").Value = 2653.6667
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 2860.375
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 2860.375
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -5856.375
# Row 122
$ws.Range("H122").Value = 1676.1666
$ws.Range("I122").Value = 2032.2222
$ws.Range("J122").Value = 1523.5714
$ws.Range("K122").Value = 6096.6666
$ws.Range("L122").Value = 4570.7142
$ws.Range("M122").Value = -3646.6666
$ws.Range("N122").Value = -9470.7142
# Row 126
$ws.Range("H126").Value = 2653.6667
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 2860.375
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 8581.125
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -13521.125
# Row 132
$ws.Range("H132").Value = 3024.2546
$ws.Range("I132").Value = 1994.5807
$ws.Range("J132").Value = 4354.25
$ws.Range("K132").Value = 5983.742099999999
$ws.Range("L132").Value = 13062.75
$ws.Range("M132").Value = -3453.742099999999
$ws.Range("N132").Value = -18122.75
# Row 134
$ws.Range("H134").Value = 2294.1
$ws.Range("I134").Value = 1167.25
$ws.Range("J134").Value = 3984.375
$ws.Range("K134").Value = 3501.75
$ws.Range("L134").Value = 11953.125
$ws.Range("M134").Value = -966.75
$ws.Range("N134").Value = -17023.125
# Row 136
$ws.Range("H136").Value = 19233526
$ws.Range("I136").Value = 1398.2142
$ws.Range("J136").Value = 41671010
$ws.Range("K136").Value = 4194.642599999999
$ws.Range("L136").Value = 125013030
$ws.Range("M136").Value = -1644.642599999999
$ws.Range("N136").Value = -125018130

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 10600
$ws.Range("I87").Value = 7960
$ws.Range("K87").Value = 23880
$ws.Range("M87").Value = -22632
# Row 90
$ws.Range("H90").Value = 10600
$ws.Range("I90").Value = 7960
$ws.Range("K90").Value = 71640
$ws.Range("M90").Value = -65400
# Row 131
$ws.Range("H131").Value = 1058.5209
$ws.Range("I131").Value = 861.5
$ws.Range("J131").Value = 1124.1945
$ws.Range("K131").Value = 2584.5
$ws.Range("L131").Value = 3372.5835
$ws.Range("M131").Value = 2455.5
$ws.Range("N131").Value = -13452.5835

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2398.8076
$ws.Range("I102").Value = 1708.9474
$ws.Range("J102").Value = 4271.2856
$ws.Range("K102").Value = 1708.9474
$ws.Range("L102").Value = 4271.2856
$ws.Range("M102").Value = -86.94740000000002
$ws.Range("N102").Value = -7515.2856
# Row 122
$ws.Range("H122").Value = 4922.1665
$ws.Range("I122").Value = 7759.8
$ws.Range("K122").Value = 23279.4
$ws.Range("M122").Value = -20829.4

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 52632790
$ws.Range("I22").Value = 250000400
$ws.Range("J22").Value = 1428.2667
$ws.Range("K22").Value = 250000400
$ws.Range("L22").Value = 1428.2667
$ws.Range("M22").Value = -250000105
$ws.Range("N22").Value = -2018.2667
# Row 27
$ws.Range("H27").Value = 52632790
$ws.Range("I27").Value = 250000400
$ws.Range("J27").Value = 1428.2667
$ws.Range("K27").Value = 250000400
$ws.Range("L27").Value = 1428.2667
$ws.Range("M27").Value = -250000293
$ws.Range("N27").Value = -1642.2667
# Row 132
$ws.Range("H132").Value = 2714.6177
$ws.Range("I132").Value = 1447.7826
$ws.Range("J132").Value = 5363.4546
$ws.Range("K132").Value = 4343.3478
$ws.Range("L132").Value = 16090.3638
$ws.Range("M132").Value = -1813.3478
$ws.Range("N132").Value = -21150.3638

$ws = $wb.Worksheets.Item("WVR")
# Row 106
$ws.Range("H106").Value = 29937.7
$ws.Range("J106").Value = 29937.7
$ws.Range("L106").Value = 29937.7
$ws.Range("N106").Value = -32461.7
# Row 108
$ws.Range("H108").Value = 29499.5
$ws.Range("J108").Value = 29499.5
$ws.Range("L108").Value = 29499.5
$ws.Range("N108").Value = -37179.5
# Row 122
$ws.Range("H122").Value = 386461.72
$ws.Range("I122").Value = 501621.1
$ws.Range("J122").Value = 2597.1667
$ws.Range("K122").Value = 1504863.3
$ws.Range("L122").Value = 7791.500100000001
$ws.Range("M122").Value = -1502413.3
$ws.Range("N122").Value = -12691.5001
# Row 132
$ws.Range("H132").Value = 273259.28
$ws.Range("I132").Value = 388977.06
$ws.Range("J132").Value = 41823.69
$ws.Range("K132").Value = 1166931.18
$ws.Range("L132").Value = 125471.07
$ws.Range("M132").Value = -1164401.18
$ws.Range("N132").Value = -130531.07
# Row 136
$ws.Range("H136").Value = 2263.25
$ws.Range("I136").Value = 1485.6364
$ws.Range("J136").Value = 3974
$ws.Range("K136").Value = 4456.9092
$ws.Range("L136").Value = 11922
$ws.Range("M136").Value = -1906.9092
